# Management Plan.xlsx edit script
# Summary of change (per commit message / diff):
#  - Rename header "01_Object" -> "01_Object Name"
#  - Insert a new "05_Tasks" column after "04_Developer"
#  - Insert a new "08_Status" column after "07_Sprint" (renumbering the
#    trailing columns 05..09 up to 06..11)
#  - Change "04_Developer" values from a person's name to an email alias,
#    and hyperlink the first occurrence (D4) to a mailto: link
#  - Populate the new "05_Tasks" / "08_Status" columns with data
#  - Append two new task rows (Customer, Product)
#  - Grow the "Tasks" table + AutoFilter + conditional formatting (data
#    bars) to cover the new columns/rows
#  - Add hyperlinks for the two new rows' "Link to Specification" cells

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planned Objects")

# ---------------------------------------------------------------------
# 1. Insert the two new columns. Doing this BEFORE touching the table
#    shifts every existing cell (values, styles, hyperlinks) along with
#    it, so none of the existing data/formatting has to be re-entered.
# ---------------------------------------------------------------------

# New column E: "05_Tasks" (pushes old E..I -> F..J)
$ws.Columns("E:E").Insert() | Out-Null
# New column H: "08_Status" (after the old "07_Sprint", now sitting in G)
$ws.Columns("H:H").Insert() | Out-Null

# ---------------------------------------------------------------------
# 2. Header row (row 3)
# ---------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value = "01_Object Name"
$ws.Cells.Item(3, 5).Value = "05_Tasks"
$ws.Cells.Item(3, 8).Value = "08_Status"
# (columns F/G/I/J/K keep their previous text: 05_Due Date, 06_Sprint,
#  07_Ready for pct, 08_Dev Comment , 09_PM Comment -- they just get
#  renumbered below)
$ws.Cells.Item(3, 6).Value = "06_Due Date"
$ws.Cells.Item(3, 7).Value = "07_Sprint"
$ws.Cells.Item(3, 9).Value = "09_Ready for pct"
$ws.Cells.Item(3, 10).Value = "10_Dev Comment "
$ws.Cells.Item(3, 11).Value = "11_PM Comment"
$ws.Rows(3).RowHeight = 43.2

# ---------------------------------------------------------------------
# 3. Developer column (D): "Sergii Razumov" -> "Sergii_Razumov@epam.com"
#    for every existing data row, then hyperlink the first cell.
# ---------------------------------------------------------------------
for ($r = 4; $r -le 8; $r++) {
    $ws.Cells.Item($r, 4).Value = "Sergii_Razumov@epam.com"
}
$ws.Hyperlinks.Add($ws.Cells.Item(4, 4), "mailto:Sergii_Razumov@epam.com") | Out-Null

# ---------------------------------------------------------------------
# 4. New "05_Tasks" column (E) values for the existing rows
# ---------------------------------------------------------------------
$ws.Cells.Item(4, 5).Value = "JR-01, JR-15, JR-16/1"
$ws.Cells.Item(5, 5).Value = "JR-02"
$ws.Cells.Item(6, 5).Value = "JR-03"
$ws.Cells.Item(7, 5).Value = "JR-04"
$ws.Cells.Item(8, 5).Value = "JR-05"

# ---------------------------------------------------------------------
# 5. New "08_Status" column (H) values for the existing rows
# ---------------------------------------------------------------------
$ws.Cells.Item(4, 8).Value = "On Hold"
$ws.Cells.Item(5, 8).Value = "On Hold"
$ws.Cells.Item(6, 8).Value = "On Hold"
$ws.Cells.Item(7, 8).Value = "To Start"
$ws.Cells.Item(8, 8).Value = "To Start"

# ---------------------------------------------------------------------
# 6. Two new rows (9 & 10), matching the style of row 4 (dates /
#    hyperlink style) by copying that row's formats first, then
#    overwriting the values.
# ---------------------------------------------------------------------
$ws.Range("A4:K4").Copy($ws.Range("A9:K9"))
$ws.Range("A4:K4").Copy($ws.Range("A10:K10"))

# Row 9: "Customer"
$ws.Cells.Item(9, 1).Value = "Customer"
$ws.Cells.Item(9, 2).Value = "Table"
$ws.Cells.Item(9, 3).Value = "https://onedrive.live.com/view.aspx?resid=43FC8CA3B17868DD%21806&id=documents&wd=target%28VS%20Code.one%7C062DAB10-4A3E-4127-8D96-1AB4A4286FC7%2FDummyTask%3A%20Add%20Internet-Sales%7CC38A5B3B-DB39-4862-8800-09D041CB42B0%2F%29"
$ws.Cells.Item(9, 4).Value = "Sergii_Razumov@epam.com"
$ws.Cells.Item(9, 5).Value = "JR-06"
$ws.Cells.Item(9, 6).Value = 44803
$ws.Cells.Item(9, 7).Value = 3
$ws.Cells.Item(9, 8).Value = "UAT"
$ws.Cells.Item(9, 9).Value = 0.95
$ws.Cells.Item(9, 10).Value = "implementing User Remarks"
$ws.Cells.Item(9, 11).Value = ""

# Row 10: "Product"
$ws.Cells.Item(10, 1).Value = "Product"
$ws.Cells.Item(10, 2).Value = "Table"
$ws.Cells.Item(10, 3).Value = "https://onedrive.live.com/view.aspx?resid=43FC8CA3B17868DD%21806&id=documents&wd=target%28VS%20Code.one%7C062DAB10-4A3E-4127-8D96-1AB4A4286FC7%2FDummyTask%3A%20Add%20Internet-Sales%7CC38A5B3B-DB39-4862-8800-09D041CB42B0%2F%29"
$ws.Cells.Item(10, 4).Value = "Sergii_Razumov@epam.com"
$ws.Cells.Item(10, 5).Value = "JR-07"
$ws.Cells.Item(10, 6).Value = 44803
$ws.Cells.Item(10, 7).Value = 3
$ws.Cells.Item(10, 8).Value = "To Start"
$ws.Cells.Item(10, 9).Value = ""
$ws.Cells.Item(10, 10).Value = "Waiting for specification"
$ws.Cells.Item(10, 11).Value = ""

# Hyperlinks for the new rows' "Link to Specification" cells
$ws.Hyperlinks.Add($ws.Cells.Item(9, 3), "https://onedrive.live.com/view.aspx?resid=43FC8CA3B17868DD%21806&id=documents&wd=target%28VS%20Code.one%7C062DAB10-4A3E-4127-8D96-1AB4A4286FC7%2FDummyTask%3A%20Add%20Internet-Sales%7CC38A5B3B-DB39-4862-8800-09D041CB42B0%2F%29") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(10, 3), "https://onedrive.live.com/view.aspx?resid=43FC8CA3B17868DD%21806&id=documents&wd=target%28VS%20Code.one%7C062DAB10-4A3E-4127-8D96-1AB4A4286FC7%2FDummyTask%3A%20Add%20Internet-Sales%7CC38A5B3B-DB39-4862-8800-09D041CB42B0%2F%29") | Out-Null

# ---------------------------------------------------------------------
# 7. Grow the "Tasks" table to cover the new columns & rows.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item("Tasks")
$lo.Resize($ws.Range("A3:K10")) | Out-Null

# ---------------------------------------------------------------------
# 8. Move the data-bar conditional formatting from the old "Ready for
#    pct" column to the new H:I range (Status + Ready for pct).
# ---------------------------------------------------------------------
$cf = $ws.Range("G4:G8").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("H4:I10")) | Out-Null

# ---------------------------------------------------------------------
# 9. Misc view bits
# ---------------------------------------------------------------------
$ws.Range("J11").Select() | Out-Null
